$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRINCIPAL")

# --- C13: the refresh re-reads this balance as a genuine number, dropping
# the stored leading zero ("0943000" text -> 943000 numeric). ---
$ws.Range("C13").Value = 943000

# --- Row 14: newly observed record appended by the automated refresh. ---
$ws.Range("A14").Value = "DF"
$ws.Range("B14").Value = "DF19110"

# C14 keeps the sheet's convention of storing this column as text, and the
# value is purely numeric digits, so force text formatting first to stop
# Excel auto-converting it to a number; then drop back to the default
# (unstyled) cell style, same as the rest of the sheet.
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "1005000"
$ws.Range("C14").Style = "Normal"

# D14 / E14 stay blank, like the rest of the sheet (a harmless no-op format
# touch materialises the otherwise-untouched cell so it is present in the
# row, matching the sheet's existing blank-cell pattern).
$ws.Range("D14").Font.Bold = $false
$ws.Range("E14").Font.Bold = $false

$ws.Range("F14").Value = "X"
$ws.Range("G14").Value = "X"
$ws.Range("H14").Value = "X - (X 01/11/25_12H) - DF"

# I14 ("01/11/25") would otherwise be auto-parsed as a date because both
# halves are valid month numbers; force text so it stays literal, matching
# the rest of the DATA_FIM column, then drop back to the default style.
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "01/11/25"
$ws.Range("I14").Style = "Normal"

$ws.Range("J14").Value = "12H"
$ws.Range("K14").Value = "19/11/25"
$ws.Range("L14").Value = "DENTRO"

# M14 (DATA_FIM_DT) stays blank for this row, same pattern as D14/E14.
$ws.Range("M14").Font.Bold = $false
